$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 5)
$ws.Range("A5").Value = "Continuing of Chapter 2 (Documentation)"

# Copy the date formatting from existing date cells so no new number
# format / style entries are created, then set the raw serial values.
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value2 = 43740

$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value2 = 43743

$ws.Range("D5").Value = "12pm"
$ws.Range("E5").Value = "5pm"

# Update selection to match the final state
$ws.Range("D6").Select()
